$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 320. This pushes the current rows
# 320-369 down to become rows 322-371, while leaving rows 318-319
# (which will be edited below) untouched.
$ws.Rows.Item(320).Insert()
$ws.Rows.Item(320).Insert()

# New row 320 receives the data that used to live in row 318
# (Murcott / Primera, Provincia de Limari).
$ws.Cells.Item(320, 1).Value = 7
$ws.Cells.Item(320, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(320, 3).Value = "Ñuble"
$ws.Cells.Item(320, 4).Value = 44469
$ws.Cells.Item(320, 5).Value = 16
$ws.Cells.Item(320, 6).Value = "Fruta"
$ws.Cells.Item(320, 7).Value = 100102
$ws.Cells.Item(320, 8).Value = "Cítricos"
$ws.Cells.Item(320, 9).Value = 100102004
$ws.Cells.Item(320, 10).Value = "Mandarina"
$ws.Cells.Item(320, 11).Value = "Murcott"
$ws.Cells.Item(320, 12).Value = "Primera"
$ws.Cells.Item(320, 13).Value = 240
$ws.Cells.Item(320, 14).Value = 6000
$ws.Cells.Item(320, 15).Value = 6500
$ws.Cells.Item(320, 16).Value = 6250
$ws.Cells.Item(320, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(320, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(320, 19).Value = 625
$ws.Cells.Item(320, 20).Value = 10

# New row 321 receives the data that used to live in row 319
# (Murcott / Segunda, Provincia de Limari).
$ws.Cells.Item(321, 1).Value = 7
$ws.Cells.Item(321, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(321, 3).Value = "Ñuble"
$ws.Cells.Item(321, 4).Value = 44469
$ws.Cells.Item(321, 5).Value = 16
$ws.Cells.Item(321, 6).Value = "Fruta"
$ws.Cells.Item(321, 7).Value = 100102
$ws.Cells.Item(321, 8).Value = "Cítricos"
$ws.Cells.Item(321, 9).Value = 100102004
$ws.Cells.Item(321, 10).Value = "Mandarina"
$ws.Cells.Item(321, 11).Value = "Murcott"
$ws.Cells.Item(321, 12).Value = "Segunda"
$ws.Cells.Item(321, 13).Value = 100
$ws.Cells.Item(321, 14).Value = 5500
$ws.Cells.Item(321, 15).Value = 5500
$ws.Cells.Item(321, 16).Value = 5500
$ws.Cells.Item(321, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(321, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(321, 19).Value = 550
$ws.Cells.Item(321, 20).Value = 10

# Row 318 is updated in place with a new record
# (Clementina / Primera, Region de O'Higgins).
$ws.Cells.Item(318, 4).Value = 45142
$ws.Cells.Item(318, 11).Value = "Clementina"
$ws.Cells.Item(318, 13).Value = 80
$ws.Cells.Item(318, 14).Value = 8000
$ws.Cells.Item(318, 15).Value = 8000
$ws.Cells.Item(318, 16).Value = 8000
$ws.Cells.Item(318, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(318, 19).Value = 800

# Row 319 is updated in place with a new record
# (Clementina / Segunda, Region de O'Higgins).
$ws.Cells.Item(319, 4).Value = 45142
$ws.Cells.Item(319, 11).Value = "Clementina"
$ws.Cells.Item(319, 13).Value = 60
$ws.Cells.Item(319, 14).Value = 6000
$ws.Cells.Item(319, 15).Value = 6000
$ws.Cells.Item(319, 16).Value = 6000
$ws.Cells.Item(319, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(319, 19).Value = 600
